# "add status in rm and pt" - replace shipment data with updated batch/status rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shipment")

$ws.Cells.Item(2, 1).Value = "TA5802HANN2"
$ws.Cells.Item(2, 2).Value = "08225SIDID"
$ws.Cells.Item(2, 3).Value = 2160
$ws.Cells.Item(2, 4).Value = 0.88

$ws.Cells.Item(3, 1).Value = "TA579F8ANN1"
$ws.Cells.Item(3, 2).Value = "08225SIDID"
$ws.Cells.Item(3, 3).Value = 2880
$ws.Cells.Item(3, 4).Value = 0.88

$ws.Cells.Item(4, 1).Value = "TA581BEANN3"
$ws.Cells.Item(4, 2).Value = "08225SIDID"
$ws.Cells.Item(4, 3).Value = 6660
$ws.Cells.Item(4, 4).Value = 0.88

$ws.Cells.Item(5, 1).Value = "TA10OF7ANV1"
$ws.Cells.Item(5, 2).Value = "08225SIDID"
$ws.Cells.Item(5, 3).Value = 1110
$ws.Cells.Item(5, 4).Value = 1.3299999999999998

$ws.Cells.Item(6, 1).Value = "TA5762NANV1"
$ws.Cells.Item(6, 2).Value = "08225SIDID"
$ws.Cells.Item(6, 3).Value = 5040
$ws.Cells.Item(6, 4).Value = 1.3299999999999998

$ws.Cells.Item(7, 1).Value = "TA10SF7ANV1"
$ws.Cells.Item(7, 2).Value = "08224SIDID"
$ws.Cells.Item(7, 3).Value = 3640
$ws.Cells.Item(7, 4).Value = 1.3299999999999998

$ws.Cells.Item(8, 1).Value = "TA10SF7ANV1"
$ws.Cells.Item(8, 2).Value = "08225SIDID"
$ws.Cells.Item(8, 3).Value = 890
$ws.Cells.Item(8, 4).Value = 1.3299999999999998

$ws.Cells.Item(9, 1).Value = "TA57797ANH2"
$ws.Cells.Item(9, 2).Value = "08225SIDID"
$ws.Cells.Item(9, 3).Value = 7400
$ws.Cells.Item(9, 4).Value = 1.3299999999999998

$ws.Cells.Item(10, 1).Value = "TA10UF7ANH2"
$ws.Cells.Item(10, 2).Value = "08225SIDID"
$ws.Cells.Item(10, 3).Value = 2090
$ws.Cells.Item(10, 4).Value = 0.73

$ws.Cells.Item(11, 1).Value = "TA10UG2ANH2"
$ws.Cells.Item(11, 2).Value = "08225SIDID"
$ws.Cells.Item(11, 3).Value = 2760
$ws.Cells.Item(11, 4).Value = 0.67

$ws.Cells.Item(12, 1).Value = "TA10VG2ANV1"
$ws.Cells.Item(12, 2).Value = "08224SIDID"
$ws.Cells.Item(12, 3).Value = 300
$ws.Cells.Item(12, 4).Value = 0.67

$ws.Cells.Item(13, 1).Value = "TA10VG2ANV1"
$ws.Cells.Item(13, 2).Value = "08225SIDID"
$ws.Cells.Item(13, 3).Value = 1810
$ws.Cells.Item(13, 4).Value = 0.67

$ws.Cells.Item(14, 1).Value = "TA10VF7ANV1"
$ws.Cells.Item(14, 2).Value = "08224SIDID"
$ws.Cells.Item(14, 3).Value = 1360
$ws.Cells.Item(14, 4).Value = 0.73

$ws.Cells.Item(15, 1).Value = "TA10VF7ANV1"
$ws.Cells.Item(15, 2).Value = "08225SIDID"
$ws.Cells.Item(15, 3).Value = 610
$ws.Cells.Item(15, 4).Value = 0.73

$ws.Cells.Item(16, 1).Value = "TA10W5RANV2"
$ws.Cells.Item(16, 2).Value = "08225SIDID"
$ws.Cells.Item(16, 3).Value = 2610
$ws.Cells.Item(16, 4).Value = 0.73

$ws.Cells.Item(17, 1).Value = "TA11J0XANV1"
$ws.Cells.Item(17, 2).Value = "08224SIDID"
$ws.Cells.Item(17, 3).Value = 950
$ws.Cells.Item(17, 4).Value = 0.67

$ws.Cells.Item(18, 1).Value = "TA11J0XANV1"
$ws.Cells.Item(18, 2).Value = "08225SIDID"
$ws.Cells.Item(18, 3).Value = 830
$ws.Cells.Item(18, 4).Value = 0.67

$ws.Cells.Item(19, 1).Value = "TA10Y0XANV1"
$ws.Cells.Item(19, 2).Value = "08225SIDID"
$ws.Cells.Item(19, 3).Value = 1250
$ws.Cells.Item(19, 4).Value = 0.67

$ws.Cells.Item(20, 1).Value = "TA11K5VANH2"
$ws.Cells.Item(20, 2).Value = "08225SIDID"
$ws.Cells.Item(20, 3).Value = 630
$ws.Cells.Item(20, 4).Value = 0.8300000000000001

$ws.Range("B7").Select()
